# Updated gitignore, ignore other changes
# Rewrite species names in column B to their .mat filenames (single-quoted),
# and update the corresponding abundance values in column C (rows 2-32).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "Akkermansia_muciniphila_ATCC_BAA_835.mat",
    "Alistipes_finegoldii_DSM_17242.mat",
    "Alistipes_putredinis_DSM_17216.mat",
    "Bacteroides_cellulosilyticus_DSM_14838.mat",
    "Bacteroides_coprophilus_DSM_18228.mat",
    "Bacteroides_fragilis_3_1_12.mat",
    "Bacteroides_oleiciplenus_YIT_12058.mat",
    "Bacteroides_ovatus_ATCC_8483.mat",
    "Bacteroides_plebeius_M12_DSM_17135.mat",
    "Bacteroides_salyersiae_WAL_10018.mat",
    "Bacteroides_thetaiotaomicron_VPI_5482.mat",
    "Bacteroides_uniformis_ATCC_8492.mat",
    "Bacteroides_vulgatus_ATCC_8482.mat",
    "Barnesiella_intestinihominis_YIT_11860.mat",
    "Bifidobacterium_animalis_lactis_AD011.mat",
    "Bilophila_wadsworthia_3_1_6.mat",
    "Escherichia_coli_O157_H7_str_Sakai.mat",
    "Eubacterium_limosum_KIST612.mat",
    "Eubacterium_ramulus_ATCC_29099.mat",
    "Flavonifractor_plautii_ATCC_29863.mat",
    "Marvinbryantia_formatexigens_I_52_DSM_14469.mat",
    "Odoribacter_splanchnicus_1651_6_DSM_20712.mat",
    "Parabacteroides_distasonis_ATCC_8503.mat",
    "Parabacteroides_johnsonii_DSM_18315.mat",
    "Paraprevotella_xylaniphila_YIT_11841.mat",
    "Parasutterella_excrementihominis_YIT_11859.mat",
    "Phascolarctobacterium_succinatutens_YIT_12067.mat",
    "Prevotella_copri_CB7_DSM_18205.mat",
    "Prevotella_stercorea_DSM_18206.mat",
    "Roseburia_inulinivorans_DSM_16841.mat",
    "Sutterella_wadsworthensis_3_1_45B.mat"
)

# Values for C2:C32 (negative-zero entries collapse to 0 - IEEE -0 == 0,
# and this runtime's xlsx writer always normalizes -0.0 to 0 on export).
$values = @(
    0,
    0,
    0,
    0.004,
    0.001,
    0.017,
    0,
    0.09,
    0,
    0.046,
    0,
    0,
    0.049,
    0,
    0,
    0,
    0.045,
    0.748,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0
)

for ($i = 0; $i -lt $names.Count; $i++) {
    $row = $i + 2
    # Leading apostrophe doubled so the literal single quote is kept as
    # text content instead of being consumed as Excel's text-prefix marker.
    $ws.Range("B$row").Value = "''" + $names[$i] + "'"
    # Clear the auto-applied "quote prefix" cell style so no style index
    # is attached to the cell (matches the unstyled source cells).
    $ws.Range("B$row").Style = "Normal"
    $ws.Range("C$row").Value = $values[$i]
}
